$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 8969.462
$ws.Range("I28").Value = 11679.223
$ws.Range("K28").Value = 11679.223
$ws.Range("M28").Value = -11194.223
$ws.Range("H52").Value = 269.55554
$ws.Range("I52").Value = 139.8
$ws.Range("J52").Value = 299.04544
$ws.Range("K52").Value = 419.4
$ws.Range("L52").Value = 897.13632
$ws.Range("M52").Value = -259.4
$ws.Range("N52").Value = -1217.13632
$ws.Range("H62").Value = 76955000
$ws.Range("I62").Value = 500010000
$ws.Range("J62").Value = 35908.547
$ws.Range("K62").Value = 500010000
$ws.Range("L62").Value = 35908.547
$ws.Range("M62").Value = -500009376
$ws.Range("N62").Value = -37156.547
$ws.Range("H64").Value = 24955.262
$ws.Range("I64").Value = 29581.777
$ws.Range("K64").Value = 29581.777
$ws.Range("M64").Value = -29333.777
$ws.Range("H65").Value = 76955000
$ws.Range("I65").Value = 500010000
$ws.Range("J65").Value = 35908.547
$ws.Range("K65").Value = 2500050000
$ws.Range("L65").Value = 179542.735
$ws.Range("M65").Value = -2500046880
$ws.Range("N65").Value = -185782.735
$ws.Range("H67").Value = 24955.262
$ws.Range("I67").Value = 29581.777
$ws.Range("K67").Value = 29581.777
$ws.Range("M67").Value = -28723.777
$ws.Range("H88").Value = 5629.8
$ws.Range("I88").Value = 4568.75
$ws.Range("J88").Value = 6015.636
$ws.Range("K88").Value = 4568.75
$ws.Range("L88").Value = 6015.636
$ws.Range("M88").Value = -4162.75
$ws.Range("N88").Value = -6827.636
$ws.Range("H91").Value = 5629.8
$ws.Range("I91").Value = 4568.75
$ws.Range("J91").Value = 6015.636
$ws.Range("K91").Value = 4568.75
$ws.Range("L91").Value = 6015.636
$ws.Range("M91").Value = -3164.75
$ws.Range("N91").Value = -8823.636
$ws.Range("H111").Value = 3791.8333
$ws.Range("I111").Value = 3791.8333
$ws.Range("K111").Value = 11375.4999
$ws.Range("M111").Value = -8308.499899999999
$ws.Range("H131").Value = 4140.316
$ws.Range("I131").Value = 1655.4546
$ws.Range("J131").Value = 7557
$ws.Range("K131").Value = 4966.3638
$ws.Range("L131").Value = 22671
$ws.Range("M131").Value = 73.63619999999992
$ws.Range("N131").Value = -32751
$ws.Range("H132").Value = 2317.537
$ws.Range("I132").Value = 2184.843
$ws.Range("J132").Value = 4573.3335
$ws.Range("K132").Value = 6554.529
$ws.Range("L132").Value = 13720.0005
$ws.Range("M132").Value = -4024.529
$ws.Range("N132").Value = -18780.0005

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1771.6786
$ws.Range("I32").Value = 1774.3148
$ws.Range("K32").Value = 1774.3148
$ws.Range("M32").Value = -1487.3148
$ws.Range("H63").Value = 6001.8
$ws.Range("J63").Value = 5003
$ws.Range("L63").Value = 5003
$ws.Range("N63").Value = -6375
$ws.Range("H66").Value = 6001.8
$ws.Range("J66").Value = 5003
$ws.Range("L66").Value = 25015
$ws.Range("N66").Value = -31879
$ws.Range("H98").Value = 134666.67
$ws.Range("J98").Value = 134666.67
$ws.Range("L98").Value = 134666.67
$ws.Range("N98").Value = -140656.67
$ws.Range("H102").Value = 19695.334
$ws.Range("I102").Value = 30764.143
$ws.Range("J102").Value = 4199
$ws.Range("K102").Value = 30764.143
$ws.Range("L102").Value = 4199
$ws.Range("M102").Value = -29142.143
$ws.Range("N102").Value = -7443
$ws.Range("H132").Value = 3134.5293
$ws.Range("I132").Value = 1461.826
$ws.Range("K132").Value = 4385.478
$ws.Range("M132").Value = -1855.478

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3166.3076
$ws.Range("I20").Value = 780
$ws.Range("J20").Value = 4657.75
$ws.Range("K20").Value = 780
$ws.Range("L20").Value = 4657.75
$ws.Range("M20").Value = -533
$ws.Range("N20").Value = -5151.75
$ws.Range("H26").Value = 45988.445
$ws.Range("I26").Value = 45988.445
$ws.Range("K26").Value = 45988.445
$ws.Range("M26").Value = -45696.445
$ws.Range("H86").Value = 11247.333
$ws.Range("I86").Value = 9134.200000000001
$ws.Range("K86").Value = 9134.200000000001
$ws.Range("M86").Value = -8011.200000000001
$ws.Range("H89").Value = 11247.333
$ws.Range("I89").Value = 9134.200000000001
$ws.Range("K89").Value = 45671
$ws.Range("M89").Value = -40055
$ws.Range("H132").Value = 71106
$ws.Range("J132").Value = 71106
$ws.Range("L132").Value = 71106
$ws.Range("N132").Value = -81226

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1478.1666
$ws.Range("I16").Value = 1478.1666
$ws.Range("K16").Value = 1478.1666
$ws.Range("M16").Value = -1191.1666
$ws.Range("H107").Value = 16543
$ws.Range("I107").Value = 24808.445
$ws.Range("J107").Value = 1665.2
$ws.Range("K107").Value = 24808.445
$ws.Range("L107").Value = 1665.2
$ws.Range("M107").Value = -22888.445
$ws.Range("N107").Value = -5505.2
$ws.Range("H113").Value = 1478.1666
$ws.Range("I113").Value = 1478.1666
$ws.Range("K113").Value = 1478.1666
$ws.Range("M113").Value = 691.8334
$ws.Range("H132").Value = 17512.385
$ws.Range("I132").Value = 1377.8948
$ws.Range("K132").Value = 4133.6844
$ws.Range("M132").Value = -1603.6844

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H110").Value = 44949.5
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()
$ws.Range("H137").Value = 2651.6
$ws.Range("I137").Value = 2847.6667
$ws.Range("J137").Value = 2357.5
$ws.Range("K137").Value = 8543.000100000001
$ws.Range("L137").Value = 7072.5
$ws.Range("M137").Value = -3443.000100000001
$ws.Range("N137").Value = -17272.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5861.5625
$ws.Range("I80").Value = 7698.8887
$ws.Range("J80").Value = 3499.2856
$ws.Range("K80").Value = 7698.8887
$ws.Range("L80").Value = 3499.2856
$ws.Range("M80").Value = -6700.8887
$ws.Range("N80").Value = -5495.2856
$ws.Range("H83").Value = 5861.5625
$ws.Range("I83").Value = 7698.8887
$ws.Range("J83").Value = 3499.2856
$ws.Range("K83").Value = 38494.4435
$ws.Range("L83").Value = 17496.428
$ws.Range("M83").Value = -33502.4435
$ws.Range("N83").Value = -27480.428
$ws.Range("H113").Value = 3832.8333
$ws.Range("J113").Value = 3799.4
$ws.Range("L113").Value = 3799.4
$ws.Range("N113").Value = -8139.4

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 28029.75
$ws.Range("I7").Value = 51516.555
$ws.Range("K7").Value = 51516.555
$ws.Range("M7").Value = -51404.555
$ws.Range("H22").Value = 2805.1
$ws.Range("I22").Value = 3242.75
$ws.Range("J22").Value = 1054.5
$ws.Range("K22").Value = 3242.75
$ws.Range("L22").Value = 1054.5
$ws.Range("M22").Value = -2947.75
$ws.Range("N22").Value = -1644.5
$ws.Range("H27").Value = 2805.1
$ws.Range("I27").Value = 3242.75
$ws.Range("J27").Value = 1054.5
$ws.Range("K27").Value = 3242.75
$ws.Range("L27").Value = 1054.5
$ws.Range("M27").Value = -3135.75
$ws.Range("N27").Value = -1268.5
$ws.Range("H46").Value = 2140.5715
$ws.Range("I46").Value = 2264.25
$ws.Range("J46").Value = 2064.4614
$ws.Range("K46").Value = 2264.25
$ws.Range("L46").Value = 2064.4614
$ws.Range("M46").Value = -2076.25
$ws.Range("N46").Value = -2440.4614
$ws.Range("H48").Value = 30041
$ws.Range("I48").Value = 30041
$ws.Range("K48").Value = 30041
$ws.Range("M48").Value = -29380
$ws.Range("H55").Value = 1422.7858
$ws.Range("I55").Value = 274
$ws.Range("J55").Value = 2571.5715
$ws.Range("K55").Value = 274
$ws.Range("L55").Value = 2571.5715
$ws.Range("M55").Value = -101
$ws.Range("N55").Value = -2917.5715
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
$ws.Range("H108").Value = 49999.5
$ws.Range("J108").Value = 49999.5
$ws.Range("L108").Value = 49999.5
$ws.Range("N108").Value = -57679.5
$ws.Range("H109").Value = 77500
$ws.Range("J109").Value = 77500
$ws.Range("L109").Value = 77500
$ws.Range("N109").Value = -80274
$ws.Range("H126").Value = 28029.75
$ws.Range("I126").Value = 51516.555
$ws.Range("K126").Value = 154549.665
$ws.Range("M126").Value = -152079.665
$ws.Range("H132").Value = 320415.38
$ws.Range("I132").Value = 553927.3
$ws.Range("K132").Value = 1661781.9
$ws.Range("M132").Value = -1659251.9

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2478.7273
$ws.Range("I96").Value = 1966.7142
$ws.Range("K96").Value = 1966.7142
$ws.Range("M96").Value = -593.7141999999999
$ws.Range("H122").Value = 27228.652
$ws.Range("I122").Value = 4804.846
$ws.Range("J122").Value = 56379.6
$ws.Range("K122").Value = 14414.538
$ws.Range("L122").Value = 169138.8
$ws.Range("M122").Value = -11964.538
$ws.Range("N122").Value = -174038.8
$ws.Range("H135").Value = 8440876
$ws.Range("J135").Value = 8440876
$ws.Range("L135").Value = 8440876
$ws.Range("N135").Value = -8451016
